# Horarios Línea 141 - actualización 03:48:17
$wb = $excel.ActiveWorkbook

$oldTimestamp = "03:18:26"
$newTimestamp = "03:48:17"

# -----------------------------------------------------------------
# Sheet 1: LP1912
# -----------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("LP1912")

$ws1.Range("A2").Value = "Última actualización: " + $newTimestamp
$ws1.Range("A3").Value = "Total filas: 8"

# existing rows 6-10 get refreshed scrap time / arrival / minutes,
# row 8 also changes line+stop name
$ws1.Cells.Item(6, 1).Value  = $newTimestamp
$ws1.Cells.Item(6, 2).Value  = "03:51"
$ws1.Cells.Item(6, 3).Value  = "14_ABASTO"
$ws1.Cells.Item(6, 4).Value  = 3
$ws1.Cells.Item(6, 5).Value  = "LP1912"

$ws1.Cells.Item(7, 1).Value  = $newTimestamp
$ws1.Cells.Item(7, 2).Value  = "04:02"
$ws1.Cells.Item(7, 3).Value  = "81_EL PELIGRO"
$ws1.Cells.Item(7, 4).Value  = 14
$ws1.Cells.Item(7, 5).Value  = "LP1912"

$ws1.Cells.Item(8, 1).Value  = $newTimestamp
$ws1.Cells.Item(8, 2).Value  = "04:47"
$ws1.Cells.Item(8, 3).Value  = "215_EL PELIGRO"
$ws1.Cells.Item(8, 4).Value  = 59
$ws1.Cells.Item(8, 5).Value  = "LP1912"

$ws1.Cells.Item(9, 1).Value  = $newTimestamp
$ws1.Cells.Item(9, 2).Value  = "04:53"
$ws1.Cells.Item(9, 3).Value  = "11_ETCHEVERRY"
$ws1.Cells.Item(9, 4).Value  = 65
$ws1.Cells.Item(9, 5).Value  = "LP1912"

$ws1.Cells.Item(10, 1).Value = $newTimestamp
$ws1.Cells.Item(10, 2).Value = "05:11"
$ws1.Cells.Item(10, 3).Value = "17_ROMERO"
$ws1.Cells.Item(10, 4).Value = 83
$ws1.Cells.Item(10, 5).Value = "LP1912"

# new rows 11-13
$ws1.Cells.Item(11, 1).Value = $newTimestamp
$ws1.Cells.Item(11, 2).Value = "05:22"
$ws1.Cells.Item(11, 3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(11, 4).Value = 94
$ws1.Cells.Item(11, 5).Value = "LP1912"

$ws1.Cells.Item(12, 1).Value = $newTimestamp
$ws1.Cells.Item(12, 2).Value = "05:32"
$ws1.Cells.Item(12, 3).Value = "81_EL PELIGRO"
$ws1.Cells.Item(12, 4).Value = 104
$ws1.Cells.Item(12, 5).Value = "LP1912"

$ws1.Cells.Item(13, 1).Value = $newTimestamp
$ws1.Cells.Item(13, 2).Value = "05:46"
$ws1.Cells.Item(13, 3).Value = "14_ABASTO"
$ws1.Cells.Item(13, 4).Value = 118
$ws1.Cells.Item(13, 5).Value = "LP1912"

# -----------------------------------------------------------------
# Sheet 2: LP1912-215
# -----------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("LP1912-215")

$ws2.Range("A2").Value = "Última actualización: " + $newTimestamp

$ws2.Cells.Item(6, 1).Value = $newTimestamp
$ws2.Cells.Item(6, 2).Value = "04:47"
$ws2.Cells.Item(6, 3).Value = "215_EL PELIGRO"
$ws2.Cells.Item(6, 4).Value = 59
$ws2.Cells.Item(6, 5).Value = "LP1912"

# -----------------------------------------------------------------
# Sheet 3: 6203-6173
# -----------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("6203-6173")

$ws3.Range("A2").Value = "Última actualización: " + $newTimestamp
